$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-09 Sunday" "2025-11-10 Monday"

Replace-Text "117÷2=58, 1" "435÷5=87, 0"
Replace-Text "268÷4=67, 0" "852÷2=426, 0"
Replace-Text "785÷3=261, 2" "137÷8=17, 1"
Replace-Text "224÷2=112, 0" "109÷3=36, 1"
Replace-Text "255÷7=36, 3" "510÷5=102, 0"

Replace-Text "709÷8=88, 5" "144÷2=72, 0"
Replace-Text "440÷5=88, 0" "953÷5=190, 3"
Replace-Text "293÷3=97, 2" "762÷4=190, 2"
Replace-Text "467÷7=66, 5" "953÷2=476, 1"
Replace-Text "337÷3=112, 1" "223÷2=111, 1"

Replace-Text "314÷3=104, 2" "389÷9=43, 2"
Replace-Text "848÷5=169, 3" "556÷7=79, 3"
Replace-Text "975÷9=108, 3" "286÷4=71, 2"
Replace-Text "429÷7=61, 2" "229÷3=76, 1"
Replace-Text "996÷8=124, 4" "403÷7=57, 4"

Replace-Text "597÷4=149, 1" "679÷8=84, 7"
Replace-Text "470÷4=117, 2" "641÷3=213, 2"
Replace-Text "598÷4=149, 2" "663÷3=221, 0"
Replace-Text "675÷4=168, 3" "370÷3=123, 1"
Replace-Text "798÷7=114, 0" "881÷9=97, 8"

Replace-Text "902÷5=180, 2" "380÷7=54, 2"
Replace-Text "747÷7=106, 5" "311÷3=103, 2"
Replace-Text "190÷4=47, 2" "748÷6=124, 4"
Replace-Text "297÷4=74, 1" "755÷3=251, 2"
Replace-Text "424÷3=141, 1" "751÷2=375, 1"
